$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained a new "2021" column (column R), mirroring the layout of the
# existing "2020" column (column Q). First copy Q2:Q6 formatting into R2:R6 so
# the new column's cell styles (borders, number formats, fonts) match column Q,
# then fill in the 2021 values.
$ws.Range("Q2:Q6").Copy() | Out-Null
$ws.Range("R2:R6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 233306
$ws.Range("R5").Value = 3.5
$ws.Range("R6").Value = 30.8

# Move the active selection, matching the saved view state of the workbook.
$ws.Range("Q15").Select() | Out-Null
